$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044646768157547
$ws.Range("D2").Value = 1.04512003900559
$ws.Range("E2").Value = 1.058249945276082
$ws.Range("F2").Value = 1.065682730545032
$ws.Range("I2").Value = 1.042725587878068
$ws.Range("J2").Value = 1.049710736182235
$ws.Range("K2").Value = 1.047889172252157
$ws.Range("L2").Value = 1.060982693753523
$ws.Range("M2").Value = 1.068395311918511
$ws.Range("N2").Value = 1.051201445831118

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045561053968665
$ws.Range("D3").Value = 1.045790769768434
$ws.Range("E3").Value = 1.059154746278503
$ws.Range("F3").Value = 1.06666796410231
$ws.Range("I3").Value = 1.042951099049845
$ws.Range("J3").Value = 1.050272536288112
$ws.Range("K3").Value = 1.048371645677239
$ws.Range("L3").Value = 1.061701260486497
$ws.Range("M3").Value = 1.069195567518689
$ws.Range("N3").Value = 1.051764043757584

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046153212857234
$ws.Range("D4").Value = 1.046225205237295
$ws.Range("E4").Value = 1.059741137817601
$ws.Range("F4").Value = 1.067306548683802
$ws.Range("I4").Value = 1.043096130744295
$ws.Range("J4").Value = 1.050635965484435
$ws.Range("K4").Value = 1.048683583522318
$ws.Range("L4").Value = 1.062166509932743
$ws.Range("M4").Value = 1.069713837830915
$ws.Range("N4").Value = 1.052127989065053

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046402288109707
$ws.Range("D5").Value = 1.046407942888794
$ws.Range("E5").Value = 1.059987876774173
$ws.Range("F5").Value = 1.067575264982005
$ws.Range("I5").Value = 1.043156888640093
$ws.Range("J5").Value = 1.050788727825055
$ws.Range("K5").Value = 1.048814660211162
$ws.Range("L5").Value = 1.062362168712369
$ws.Range("M5").Value = 1.069931825469625
$ws.Range("N5").Value = 1.052280968345717

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04644411662697
$ws.Range("D6").Value = 1.046438631224253
$ws.Range("E6").Value = 1.060029318186985
$ws.Range("F6").Value = 1.067620398562373
$ws.Range("I6").Value = 1.043167077625985
$ws.Range("J6").Value = 1.050814375898684
$ws.Range("K6").Value = 1.048836664900334
$ws.Range("L6").Value = 1.062395024629766
$ws.Range("M6").Value = 1.069968432779185
$ws.Range("N6").Value = 1.052306652842554

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046156540497086
$ws.Range("D7").Value = 1.04622764659176
$ws.Range("E7").Value = 1.05974443389315
$ws.Range("F7").Value = 1.067310138282614
$ws.Range("I7").Value = 1.043096943433611
$ws.Range("J7").Value = 1.050638006795108
$ws.Range("K7").Value = 1.048685335220243
$ws.Range("L7").Value = 1.06216912406959
$ws.Range("M7").Value = 1.069716750174351
$ws.Range("N7").Value = 1.05213003327462

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044955639740178
$ws.Range("D8").Value = 1.045346626160113
$ws.Range("E8").Value = 1.0585555346595
$ws.Range("F8").Value = 1.066015471885082
$ws.Range("I8").Value = 1.042801984185942
$ws.Range("J8").Value = 1.049900617808875
$ws.Range("K8").Value = 1.048052278545555
$ws.Range("L8").Value = 1.061225476236725
$ws.Range("M8").Value = 1.068665667807966
$ws.Range("N8").Value = 1.051391597111448

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042843800406059
$ws.Range("D9").Value = 1.043797495732005
$ws.Range("E9").Value = 1.056467688883447
$ws.Range("F9").Value = 1.063742382738115
$ws.Range("I9").Value = 1.042275444148931
$ws.Range("J9").Value = 1.048600576734488
$ws.Range("K9").Value = 1.046934845139023
$ws.Range("L9").Value = 1.059564912869203
$ws.Range("M9").Value = 1.066817040105329
$ws.Range("N9").Value = 1.050089709829618

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041438868508522
$ws.Range("D10").Value = 1.042767073296283
$ws.Range("E10").Value = 1.055080682115223
$ws.Range("F10").Value = 1.062232641172964
$ws.Range("I10").Value = 1.041919890195743
$ws.Range("J10").Value = 1.047733494282127
$ws.Range("K10").Value = 1.046188666878364
$ws.Range("L10").Value = 1.058459464782678
$ws.Range("M10").Value = 1.065587059431108
$ws.Range("N10").Value = 1.04922139602072

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040831234074632
$ws.Range("D11").Value = 1.042321459872372
$ws.Range("E11").Value = 1.054481270657564
$ws.Range("F11").Value = 1.061580264556375
$ws.Range("I11").Value = 1.041764864692881
$ws.Range("J11").Value = 1.047357958464149
$ws.Range("K11").Value = 1.045865285841159
$ws.Range("L11").Value = 1.057981186580889
$ws.Range("M11").Value = 1.0650550581359
$ws.Range("N11").Value = 1.048845326898806

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040605639425992
$ws.Range("D12").Value = 1.042156025722083
$ws.Range("E12").Value = 1.054258799996697
$ws.Range("F12").Value = 1.061338147281941
$ws.Range("I12").Value = 1.041707121367003
$ws.Range("J12").Value = 1.047218456201542
$ws.Range("K12").Value = 1.04574512649691
$ws.Range("L12").Value = 1.057803592249406
$ws.Range("M12").Value = 1.064857538816449
$ws.Range("N12").Value = 1.048705626526984

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040654025346233
$ws.Range("D13").Value = 1.042191507970782
$ws.Range("E13").Value = 1.054306512656142
$ws.Range("F13").Value = 1.061390072984828
$ws.Range("I13").Value = 1.041719514735746
$ws.Range("J13").Value = 1.047248380416672
$ws.Range("K13").Value = 1.045770902934108
$ws.Range("L13").Value = 1.057841684125465
$ws.Range("M13").Value = 1.064899903304209
$ws.Range("N13").Value = 1.048735593237932

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040812584144339
$ws.Range("D14").Value = 1.042307783260684
$ws.Range("E14").Value = 1.054462877535326
$ws.Range("F14").Value = 1.061560246888101
$ws.Range("I14").Value = 1.041760094871338
$ws.Range("J14").Value = 1.047346427393263
$ws.Range("K14").Value = 1.045855354267669
$ws.Range("L14").Value = 1.057966505357686
$ws.Range("M14").Value = 1.065038729291339
$ws.Range("N14").Value = 1.048833779452476

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040910291775539
$ws.Range("D15").Value = 1.042379435811599
$ws.Range("E15").Value = 1.054559242653558
$ws.Range("F15").Value = 1.061665123779914
$ws.Range("I15").Value = 1.04178507645348
$ws.Range("J15").Value = 1.047406835872113
$ws.Range("K15").Value = 1.045907382087613
$ws.Range("L15").Value = 1.058043419756506
$ws.Range("M15").Value = 1.065124276480367
$ws.Range("N15").Value = 1.048894273718295

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041479210265826
$ws.Range("D16").Value = 1.042796659245448
$ws.Range("E16").Value = 1.055120487885492
$ws.Range("F16").Value = 1.062275965853997
$ws.Range("I16").Value = 1.04193015626674
$ws.Range("J16").Value = 1.047758415684998
$ws.Range("K16").Value = 1.046210122777383
$ws.Range("L16").Value = 1.058491214809391
$ws.Range("M16").Value = 1.065622379108809
$ws.Range("N16").Value = 1.049246352814841

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041836268768814
$ws.Range("D17").Value = 1.043058525109678
$ws.Range("E17").Value = 1.055472856870542
$ws.Range("F17").Value = 1.062659493762119
$ws.Range("I17").Value = 1.042020875373508
$ws.Range("J17").Value = 1.047978930784152
$ws.Range("K17").Value = 1.046399949491417
$ws.Range("L17").Value = 1.058772209555629
$ws.Range("M17").Value = 1.065934984148245
$ws.Range("N17").Value = 1.049467181070729

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042044603416076
$ws.Range("D18").Value = 1.043211321580318
$ws.Range("E18").Value = 1.05567850068746
$ws.Range("F18").Value = 1.062883329402198
$ws.Range("I18").Value = 1.042073687181805
$ws.Range("J18").Value = 1.048107545397821
$ws.Range("K18").Value = 1.046510645031145
$ws.Range("L18").Value = 1.058936146405316
$ws.Range("M18").Value = 1.066117378072728
$ws.Range("N18").Value = 1.049595978331902

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04211565172574
$ws.Range("D19").Value = 1.043263430435792
$ws.Range("E19").Value = 1.055748639029763
$ws.Range("F19").Value = 1.062959673629375
$ws.Range("I19").Value = 1.042091677137479
$ws.Range("J19").Value = 1.048151398251217
$ws.Range("K19").Value = 1.046548384711181
$ws.Range("L19").Value = 1.058992050937052
$ws.Range("M19").Value = 1.066179579266064
$ws.Range("N19").Value = 1.04963989346138

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041797952681059
$ws.Range("D20").Value = 1.043030423728752
$ws.Range("E20").Value = 1.055435039320067
$ws.Range("F20").Value = 1.062618331359517
$ws.Range("I20").Value = 1.042011152733294
$ws.Range("J20").Value = 1.047955272432938
$ws.Range("K20").Value = 1.046379585684857
$ws.Range("L20").Value = 1.058742057604971
$ws.Range("M20").Value = 1.065901438712898
$ws.Range("N20").Value = 1.049443489121943

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.04076588954453
$ws.Range("D21").Value = 1.042273540686078
$ws.Range("E21").Value = 1.054416827065979
$ws.Range("F21").Value = 1.06151012925996
$ws.Range("I21").Value = 1.041748149443875
$ws.Range("J21").Value = 1.047317555312987
$ws.Range("K21").Value = 1.04583048658508
$ws.Range("L21").Value = 1.057929746964248
$ws.Range("M21").Value = 1.064997846013547
$ws.Range("N21").Value = 1.048804866370534

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040117614943538
$ws.Range("D22").Value = 1.041798159544283
$ws.Range("E22").Value = 1.05377766394613
$ws.Range("F22").Value = 1.060814542189147
$ws.Range("I22").Value = 1.041581863679684
$ws.Range("J22").Value = 1.046916530978973
$ws.Range("K22").Value = 1.045485008137071
$ws.Range("L22").Value = 1.057419359788415
$ws.Range("M22").Value = 1.064430240717779
$ws.Range("N22").Value = 1.048403272535968

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040461218010945
$ws.Range("D23").Value = 1.042050120087814
$ws.Range("E23").Value = 1.054116398554858
$ws.Range("F23").Value = 1.061183173395749
$ws.Range("I23").Value = 1.041670102420178
$ws.Range("J23").Value = 1.0471291275165
$ws.Range("K23").Value = 1.045668175036632
$ws.Range("L23").Value = 1.057689892613327
$ws.Range("M23").Value = 1.064731089385942
$ws.Range("N23").Value = 1.048616170984963

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041815265860375
$ws.Range("D24").Value = 1.043043121364823
$ws.Range("E24").Value = 1.055452127094444
$ws.Range("F24").Value = 1.062636930474976
$ws.Range("I24").Value = 1.042015546294553
$ws.Range("J24").Value = 1.047965962648657
$ws.Range("K24").Value = 1.046388787296548
$ws.Range("L24").Value = 1.058755681858256
$ws.Range("M24").Value = 1.065916596276134
$ws.Range("N24").Value = 1.049454194518994

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043389244676843
$ws.Range("D25").Value = 1.044197578448414
$ws.Range("E25").Value = 1.057006591359434
$ws.Range("F25").Value = 1.064329040959822
$ws.Range("I25").Value = 1.042412367831178
$ws.Range("J25").Value = 1.048936741400266
$ws.Range("K25").Value = 1.047223948306188
$ws.Range("L25").Value = 1.059993932567644
$ws.Range("M25").Value = 1.06729453053028
$ws.Range("N25").Value = 1.05042635188778
